$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("425:425").Insert()

$ws.Range("A425").Value = 4
$ws.Range("B425").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C425").Value = "Los Lagos"
$ws.Range("D425").Value = 45173
$ws.Range("E425").Value = 10
$ws.Range("F425").Value = 100112003
$ws.Range("G425").Value = "Ajo"
$ws.Range("H425").Value = "Chino"
$ws.Range("I425").Value = "Primera"
$ws.Range("J425").Value = 80
$ws.Range("K425").Value = 23000
$ws.Range("L425").Value = 23000
$ws.Range("M425").Value = 23000
$ws.Range("N425").Value = "$/caja 10 kilos"
$ws.Range("O425").Value = "China"
$ws.Range("P425").Value = 2300
$ws.Range("Q425").Value = 10
$ws.Range("R425").Value = "Hortaliza"
